$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4528
$ws.Range("C2").Value = 5251

$ws.Range("B3").Value = 2983
$ws.Range("C3").Value = 3378

$ws.Range("B4").Value = 1545
$ws.Range("C4").Value = 1873

$ws.Range("B5").Value = 1022
$ws.Range("C5").Value = 1134

$ws.Range("B6").Value = 304
$ws.Range("C6").Value = 256

$ws.Range("B7").Value = 1240
$ws.Range("C7").Value = 1356

$ws.Range("B8").Value = 45
$ws.Range("C8").Value = 47

$ws.Range("B9").Value = 1070
$ws.Range("C9").Value = 1197

$ws.Range("B10").Value = 94
$ws.Range("C10").Value = 371

$ws.Range("F26").Select()
